$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @(
    "81-76=",
    "39+59=",
    "13-5=",
    "19+17=",
    "81-4=",
    "52-15=",
    "22-6=",
    "80-49=",
    "83-48=",
    "15+79=",
    "27+34=",
    "51-24=",
    "47+5=",
    "68+6=",
    "5+38=",
    "79+16=",
    "90-36=",
    "42-4=",
    "15+28=",
    "29+19=",
    "7+54=",
    "90-41=",
    "57-49=",
    "61-58=",
    "28+23=",
    "26+66=",
    "31-14=",
    "28+17=",
    "80-43=",
    "24-16=",
    "90-87=",
    "36+57=",
    "64-19=",
    "73-69=",
    "39+22=",
    "37+49=",
    "68+28=",
    "56-17=",
    "19+69=",
    "70-17=",
    "66+7=",
    "61-46=",
    "91-45=",
    "6+49=",
    "73-27=",
    "70-48=",
    "90-47=",
    "65-56=",
    "19+9=",
    "67-49=",
    "6+67=",
    "39+26=",
    "71-46=",
    "90-89=",
    "68-9=",
    "44-25=",
    "74-38=",
    "26+17=",
    "19+52=",
    "29+37=",
    "91-88=",
    "57+4=",
    "55-49=",
    "57+25=",
    "13-7=",
    "71-43=",
    "15+39=",
    "15+46=",
    "58+5=",
    "72-37=",
    "70-39=",
    "17+38=",
    "84-67=",
    "78+17=",
    "5+48=",
    "45+6=",
    "46-17=",
    "23-7=",
    "29+65=",
    "90-23=",
    "58+26=",
    "43-14=",
    "20-8=",
    "62+19=",
    "57+37=",
    "73-7=",
    "67-48=",
    "94-38=",
    "45+18=",
    "90-14=",
    "47+18=",
    "76-29=",
    "79+19=",
    "15+28=",
    "15+76=",
    "94-39=",
    "70-17=",
    "75-68=",
    "72-67=",
    "36+47="
)

$cols = $t.Columns.Count
$index = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $t.Cell($r, $c).Range.Text = $newValues[$index]
        $index = $index + 1
    }
}

Write-Host "Updated $index cells"
